$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.973.08"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "1.649.60"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = "  +0.58%  "
$ws.Range("D5").Value = "'216.01"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("D6").Value = "'0.5112"
$ws.Range("E6").Value = "  +1.99%  "
$ws.Range("D7").Value = "'1.006"
$ws.Range("E7").Value = "  +0.40%  "
$ws.Range("D8").Value = "'0.2587"
$ws.Range("E8").Value = "  +0.80%  "
$ws.Range("D9").Value = "'0.06433"
$ws.Range("E9").Value = "  +0.54%  "
$ws.Range("D10").Value = "'19.75"
$ws.Range("E10").Value = "  +0.80%  "
$ws.Range("D11").Value = "'0.07782"
$ws.Range("E11").Value = "  +1.36%  "
$ws.Range("D12").Value = "'4.326"
$ws.Range("E12").Value = "  +2.07%  "
$ws.Range("D13").Value = "1.646.74"
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("D14").Value = "'0.5489"
$ws.Range("E14").Value = "  +1.16%  "
$ws.Range("D15").Value = "0.0₅7899"
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("D16").Value = "'65.00"
$ws.Range("E16").Value = "  +2.55%  "
$ws.Range("D17").Value = "26.010.18"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("D18").Value = "'1.006"
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("D19").Value = "'199.24"
$ws.Range("E19").Value = "  -1.35%  "
$ws.Range("D20").Value = "'4.459"
$ws.Range("E20").Value = "  +3.11%  "
$ws.Range("D21").Value = "'10.06"
$ws.Range("E21").Value = "  +1.36%  "
$ws.Range("D22").Value = "'6.076"
$ws.Range("E22").Value = "  +1.84%  "
$ws.Range("D23").Value = "'1.007"
$ws.Range("E23").Value = "  +0.43%  "
$ws.Range("D24").Value = "'1.865"
$ws.Range("E24").Value = "  -2.46%  "
$ws.Range("D25").Value = "'140.55"
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("D26").Value = "'0.1153"
$ws.Range("E26").Value = "  +1.32%  "
$ws.Range("D27").Value = "'6.919"
$ws.Range("E27").Value = "  +3.40%  "
$ws.Range("D28").Value = "'15.79"
$ws.Range("E28").Value = "  +0.86%  "
$ws.Range("D29").Value = "'1.242"
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("D30").Value = "'0.05039"
$ws.Range("E30").Value = "  +1.27%  "
$ws.Range("D31").Value = "'3.294"
$ws.Range("E31").Value = "  +1.26%  "
$ws.Range("D32").Value = "'3.212"
$ws.Range("E32").Value = "  +1.23%  "
$ws.Range("D33").Value = "'1.546"
$ws.Range("E33").Value = "  +0.60%  "
$ws.Range("D34").Value = "'2.367"
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("D35").Value = "'0.8969"
$ws.Range("E35").Value = "  +0.56%  "
$ws.Range("D36").Value = "'2.591"
$ws.Range("E36").Value = "  -0.91%  "
$ws.Range("D37").Value = "1.139.58"
$ws.Range("E37").Value = "  -2.24%  "
$ws.Range("D38").Value = "'0.5560"
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("D39").Value = "'0.01566"
$ws.Range("E39").Value = "  +0.72%  "
$ws.Range("D40").Value = "'1.007"
$ws.Range("E40").Value = "  +0.49%  "
$ws.Range("D41").Value = "'5.675"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").Value = "'0.8181"
$ws.Range("E42").Value = "  +1.50%  "
$ws.Range("D43").Value = "'99.96"
$ws.Range("E43").Value = "  +0.81%  "
$ws.Range("E44").Value = "  +8.65%  "
$ws.Range("D45").Value = "1.783.43"
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("D46").Value = "'0.4537"
$ws.Range("E46").Value = "  +0.52%  "
$ws.Range("D47").Value = "'55.44"
$ws.Range("E47").Value = "  +1.37%  "
$ws.Range("D48").Value = "'1.006"
$ws.Range("E48").Value = "  +0.33%  "
$ws.Range("D49").Value = "'0.05099"
$ws.Range("E49").Value = "  +0.44%  "
$ws.Range("D50").Value = "'0.09592"
$ws.Range("E50").Value = "  +3.61%  "
$ws.Range("D51").Value = "'1.004"
$ws.Range("E51").Value = "  +0.13%  "
